$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add Patrick Malzahn's result row (row 33)
$ws.Range("B33").Value = "SMF"
$ws.Range("C33").Value = "PATRICK MALZAHN"
$ws.Range("D33").Value = " 1:04.626"

# Update sheet view/selection state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("D35").Select()
